$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.472488522529602
$ws.Range("B1").Value = 1.967365622520447
$ws.Range("C1").Value = 2.221557140350342
$ws.Range("D1").Value = 2.521106481552124
$ws.Range("E1").Value = 3.101007699966431
